$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Cd86"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 81.47348266666667
$ws.Range("H2").Value = 244.420448
$ws.Range("I2").Value = 0.3594530042390097
$ws.Range("J2").Value = 0.3594530042390096
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.430875666666666
$ws.Range("N2").Value = 13.292627
$ws.Range("O2").Value = 0.3569860433688579
$ws.Range("P2").Value = 0.3569860433688579
$ws.Range("Q2").Value = 360.9988718263218
$ws.Range("R2").Value = 3248.989846436896
$ws.Range("S2").Value = 0.1283197057603334
$ws.Range("T2").Value = 0.1283197057603333

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Cd86"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 81.47348266666667
$ws.Range("H3").Value = 244.420448
$ws.Range("I3").Value = 0.3594530042390097
$ws.Range("J3").Value = 0.3594530042390096
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.025118333333332
$ws.Range("N3").Value = 21.075355
$ws.Range("O3").Value = 0.5659985489733577
$ws.Range("P3").Value = 0.5659985489733576
$ws.Range("Q3").Value = 572.3608567621155
$ws.Range("R3").Value = 5151.247710859039
$ws.Range("S3").Value = 0.2034498788233937
$ws.Range("T3").Value = 0.2034498788233936

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Cd86"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 81.47348266666667
$ws.Range("H4").Value = 244.420448
$ws.Range("I4").Value = 0.3594530042390097
$ws.Range("J4").Value = 0.3594530042390096
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.950574
$ws.Range("N4").Value = 2.851722
$ws.Range("O4").Value = 0.0765856857013987
$ws.Range("P4").Value = 0.07658568570139869
$ws.Range("Q4").Value = 77.446574312384
$ws.Range("R4").Value = 697.019168811456
$ws.Range("S4").Value = 0.02752895480707233
$ws.Range("T4").Value = 0.02752895480707232

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Cd86"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 81.47348266666667
$ws.Range("H5").Value = 244.420448
$ws.Range("I5").Value = 0.3594530042390097
$ws.Range("J5").Value = 0.3594530042390096
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.005333666666666667
$ws.Range("N5").Value = 0.016001
$ws.Range("O5").Value = 0.0004297219563856788
$ws.Range("P5").Value = 0.0004297219563856788
$ws.Range("Q5").Value = 0.4345523987164445
$ws.Range("R5").Value = 3.910971588448
$ws.Range("S5").Value = 0.000154464848210297
$ws.Range("T5").Value = 0.0001544648482102969

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cd86"
$ws.Range("C6").Value = "Cd28"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 123.2241746666667
$ws.Range("H6").Value = 369.672524
$ws.Range("I6").Value = 0.5436529571225457
$ws.Range("J6").Value = 0.5436529571225457
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.430875666666666
$ws.Range("N6").Value = 13.292627
$ws.Range("O6").Value = 0.3569860433688579
$ws.Range("P6").Value = 0.3569860433688579
$ws.Range("Q6").Value = 545.9909970756164
$ws.Range("R6").Value = 4913.918973680548
$ws.Range("S6").Value = 0.194076518128957
$ws.Range("T6").Value = 0.1940765181289569

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd86"
$ws.Range("C7").Value = "Cd28"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 123.2241746666667
$ws.Range("H7").Value = 369.672524
$ws.Range("I7").Value = 0.5436529571225457
$ws.Range("J7").Value = 0.5436529571225457
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.025118333333332
$ws.Range("N7").Value = 21.075355
$ws.Range("O7").Value = 0.5659985489733577
$ws.Range("P7").Value = 0.5659985489733576
$ws.Range("Q7").Value = 865.6644085606688
$ws.Range("R7").Value = 7790.97967704602
$ws.Range("S7").Value = 0.3077067848764359
$ws.Range("T7").Value = 0.3077067848764358

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd86"
$ws.Range("C8").Value = "Cd28"
$ws.Range("D8").Value = "Neutro"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 123.2241746666667
$ws.Range("H8").Value = 369.672524
$ws.Range("I8").Value = 0.5436529571225457
$ws.Range("J8").Value = 0.5436529571225457
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.950574
$ws.Range("N8").Value = 2.851722
$ws.Range("O8").Value = 0.0765856857013987
$ws.Range("P8").Value = 0.07658568570139869
$ws.Range("Q8").Value = 117.133696609592
$ws.Range("R8").Value = 1054.203269486328
$ws.Range("S8").Value = 0.04163603450482327
$ws.Range("T8").Value = 0.04163603450482326

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd86"
$ws.Range("C9").Value = "Cd28"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 123.2241746666667
$ws.Range("H9").Value = 369.672524
$ws.Range("I9").Value = 0.5436529571225457
$ws.Range("J9").Value = 0.5436529571225457
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.005333666666666667
$ws.Range("N9").Value = 0.016001
$ws.Range("O9").Value = 0.0004297219563856788
$ws.Range("P9").Value = 0.0004297219563856788
$ws.Range("Q9").Value = 0.6572366729471112
$ws.Range("R9").Value = 5.915130056524001
$ws.Range("S9").Value = 0.0002336196123295599
$ws.Range("T9").Value = 0.0002336196123295599

# Row 10
$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Cd86"
$ws.Range("C10").Value = "Cd28"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 21.96196633333333
$ws.Range("H10").Value = 65.88589899999999
$ws.Range("I10").Value = 0.0968940386384447
$ws.Range("J10").Value = 0.09689403863844469
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.430875666666666
$ws.Range("N10").Value = 13.292627
$ws.Range("O10").Value = 0.3569860433688579
$ws.Range("P10").Value = 0.3569860433688579
$ws.Range("Q10").Value = 97.31074221851921
$ws.Range("R10").Value = 875.7966799666729
$ws.Range("S10").Value = 0.03458981947956762
$ws.Range("T10").Value = 0.0345898194795676

# Row 11
$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Cd86"
$ws.Range("C11").Value = "Cd28"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 21.96196633333333
$ws.Range("H11").Value = 65.88589899999999
$ws.Range("I11").Value = 0.0968940386384447
$ws.Range("J11").Value = 0.09689403863844469
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.025118333333332
$ws.Range("N11").Value = 21.075355
$ws.Range("O11").Value = 0.5659985489733577
$ws.Range("P11").Value = 0.5659985489733576
$ws.Range("Q11").Value = 154.2854123243494
$ws.Range("R11").Value = 1388.568710919145
$ws.Range("S11").Value = 0.05484188527352815
$ws.Range("T11").Value = 0.05484188527352814

# Row 12
$ws.Range("A12").Value = "Neutro"
$ws.Range("B12").Value = "Cd86"
$ws.Range("C12").Value = "Cd28"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 21.96196633333333
$ws.Range("H12").Value = 65.88589899999999
$ws.Range("I12").Value = 0.0968940386384447
$ws.Range("J12").Value = 0.09689403863844469
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.950574
$ws.Range("N12").Value = 2.851722
$ws.Range("O12").Value = 0.0765856857013987
$ws.Range("P12").Value = 0.07658568570139869
$ws.Range("Q12").Value = 20.876474185342
$ws.Range("R12").Value = 187.888267668078
$ws.Range("S12").Value = 0.007420696389503108
$ws.Range("T12").Value = 0.007420696389503105

# Row 13
$ws.Range("A13").Value = "Neutro"
$ws.Range("B13").Value = "Cd86"
$ws.Range("C13").Value = "Cd28"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 21.96196633333333
$ws.Range("H13").Value = 65.88589899999999
$ws.Range("I13").Value = 0.0968940386384447
$ws.Range("J13").Value = 0.09689403863844469
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.005333666666666667
$ws.Range("N13").Value = 0.016001
$ws.Range("O13").Value = 0.0004297219563856788
$ws.Range("P13").Value = 0.0004297219563856788
$ws.Range("Q13").Value = 0.1171378077665556
$ws.Range("R13").Value = 1.054240269899
$ws.Range("S13").Value = 0.00004163749584582201
$ws.Range("T13").Value = 0.000041637495845822
